$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the formatting from
# the existing header cell H1 so they match the other headers.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill data rows 2-29: column I = 1, column J = same value as column H.
for ($r = 2; $r -le 29; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}

# Row 30 is special-cased with explicit values.
$ws.Cells.Item(30, 9).Value = 4
$ws.Cells.Item(30, 10).Value = 5
